# Apply updated dSF (column F) values for the kershaw_clayton workbook.
# Row 6 (F6) is intentionally left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -1
    3  = -3
    4  = 5
    5  = 2
    7  = 6
    8  = -3
    9  = -2
    10 = 11
    11 = -1
    12 = -5
    13 = 7
    14 = -1
    15 = 3
    16 = 13
    17 = -4
    18 = 2
    19 = -2
    20 = 1
    21 = 1
    22 = 3
    23 = -1
    24 = -7
    25 = -5
    26 = -1
    27 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
